{"js": "// tdf#112118 fixture edit: move the \"_GoBack\" bookmark from the last\n// paragraph (which currently only holds the bookmark) to the very\n// first paragraph, replacing that first paragraph's page-break run.\n// The paragraph that used to hold the bookmark becomes empty.\n\nconst OOXML_NS =\n  'xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"';\nconst WORD_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPackage(bodyInnerXml) {\n  return (\n    '<pkg:package ' + OOXML_NS + '>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document ' + WORD_NS + '>' +\n    '<w:body>' + bodyInnerXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Empty out the paragraph that used to carry the \"_GoBack\" bookmark.\nconst emptyParagraphOoxml = wrapPackage(\"<w:p/>\");\nlastParagraph.insertOoxml(emptyParagraphOoxml, Word.InsertLocation.replace);\n\n// Replace the first paragraph's page-break run with the bookmark pair.\nconst bookmarkParagraphOoxml = wrapPackage(\n  '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n);\nfirstParagraph.insertOoxml(bookmarkParagraphOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# tdf#112118 fixture edit: move the \"_GoBack\" bookmark from the last\n# paragraph (which currently only holds the bookmark) to the very\n# first paragraph, replacing that first paragraph's page-break run.\n# The paragraph that used to hold the bookmark becomes empty.\n\n$d = $word.ActiveDocument\n\n# Explicitly remove the existing \"_GoBack\" bookmark first (Word only\n# ever keeps one bookmark per name, but do this defensively so the\n# script does not depend on that dedupe behaviour).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Remove the page-break run from the first paragraph, keeping the\n# paragraph mark itself intact.\n$firstParagraph = $d.Paragraphs(1).Range\n$runOnly = $d.Range($firstParagraph.Start, $firstParagraph.End - 1)\n$runOnly.Delete()\n\n# Re-create the \"_GoBack\" bookmark collapsed at the start of the first\n# paragraph (where the page break used to be).\n$target = $d.Paragraphs(1).Range\n$target.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $target)\n"}
